$d = $word.ActiveDocument

# Update the date/weekday heading paragraph.
$d.Paragraphs.Item(1).Range.Text = "2025-12-04 Thursday"

# Update the practice-table answers (three-digit x one-digit multiplication).
# Only the populated rows (1, 5, 10, 15, 20) contain data; the rest are blank
# spacer rows. Addressing cells directly by (row, column) avoids any
# ambiguity from the fact that some new values coincide with other old
# values elsewhere in the table (e.g. "226x7=1582" is both an old value in
# row 5 and a new value being written into row 20).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "829×3=2487"
$t.Cell(1, 2).Range.Text  = "115×7=805"
$t.Cell(1, 3).Range.Text  = "986×6=5916"
$t.Cell(1, 4).Range.Text  = "415×8=3320"
$t.Cell(1, 5).Range.Text  = "825×2=1650"

$t.Cell(5, 1).Range.Text  = "349×9=3141"
$t.Cell(5, 2).Range.Text  = "741×7=5187"
$t.Cell(5, 3).Range.Text  = "690×5=3450"
$t.Cell(5, 4).Range.Text  = "226×7=1582"
$t.Cell(5, 5).Range.Text  = "429×5=2145"

$t.Cell(10, 1).Range.Text = "535×5=2675"
$t.Cell(10, 2).Range.Text = "665×3=1995"
$t.Cell(10, 3).Range.Text = "353×8=2824"
$t.Cell(10, 4).Range.Text = "679×7=4753"
$t.Cell(10, 5).Range.Text = "900×3=2700"

$t.Cell(15, 1).Range.Text = "270×4=1080"
$t.Cell(15, 2).Range.Text = "467×5=2335"
$t.Cell(15, 3).Range.Text = "353×4=1412"
$t.Cell(15, 4).Range.Text = "917×6=5502"
$t.Cell(15, 5).Range.Text = "663×3=1989"

$t.Cell(20, 1).Range.Text = "999×7=6993"
$t.Cell(20, 2).Range.Text = "715×5=3575"
$t.Cell(20, 3).Range.Text = "569×4=2276"
$t.Cell(20, 4).Range.Text = "519×3=1557"
$t.Cell(20, 5).Range.Text = "468×5=2340"
